# Insert a new abstract paragraph immediately before the existing
# "Abstract here!" placeholder paragraph.
#
# The new paragraph is built from its exact target WordprocessingML markup
# (matching the runs/formatting captured in the authoritative diff -
# Calibri/shaded runs copy-pasted from the web for the first part of the
# abstract, then plain Calibri/Times-New-Roman-complex-script runs for the
# JPE sentence) and is inserted via Range.InsertXML so the run/proofErr
# structure is reproduced exactly.

$d = $word.ActiveDocument

# Locate the placeholder paragraph ("Abstract here!") and open a new,
# empty paragraph immediately above it.
$placeholder = $d.Paragraphs.First
$placeholder.Range.InsertParagraphBefore() | Out-Null

# The freshly created paragraph is now the first paragraph in the body;
# fill it in with the full abstract markup.
$abstractRange = $d.Paragraphs.First.Range

$abstractXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:textAlignment w:val="baseline"/>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">California Department of Water Resources (DWR) currently operates multiple eight-foot rotary screw traps (RSTs) in both the Low Flow Channel (LFC) and High Flow Channel (HFC) of the upper 21.5 river miles of the Lower Feather River. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve">The LFC extends from the Fish Barrier Dam at river mile (RM) 67.2 to the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Thermalito</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> Outlet (RM 59). The HFC extends from the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>Thermalito</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> Outlet to the confluence with</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t>the Sacramento River.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
          <w:color w:val="000000"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        </w:rPr>
        <w:t xml:space="preserve"> This survey has been ongoing since 1997 with the objective of documenting general salmonid emigration attributes, such as timing, abundance and composition of salmon and investigating the influence of environmental factors thought to initiate emigration, such as flow, turbidity and water temperature. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t>Data from this monitoring will</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t xml:space="preserve"> also</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t xml:space="preserve"> be used to inform the development of a juvenile production estimate (JPE) for spring-run Chinook salmon in the Sacramento River Watershed as required by Incidental Take Permit No. 2081-2019-006-00 issued by C</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t>A Department of Fish and Wildlife (CDFW)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t xml:space="preserve"> to </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t>DWR</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/>
        </w:rPr>
        <w:t xml:space="preserve"> for the long-term operation of the State Water Project.</w:t>
      </w:r>
    </w:p>
'@

$abstractRange.InsertXML($abstractXml)

Write-Output "Inserted abstract paragraph before the placeholder paragraph."
